# "nova rodada 1 a 9" — refresh the computed metrics (columns B:I, rows 1-5)
# for another round of the qtd-criteria experiment, restyle the now much
# larger C1 value in scientific notation, and reposition/resize the first
# scatter chart on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- updated data values --------------------------------------------------
$values = @{
    "B1" = 0.0209865996250529;  "C1" = 64048999563.880798;  "D1" = 0.855151515151515
    "E1" = 0.204846691761921;   "F1" = 0.80999999999999905; "G1" = 0.268700576850888
    "H1" = 0.96464962126200904; "I1" = 0.049992985046292002

    "B2" = 0.022931774166958199; "C2" = 0.016146803710260602; "D2" = 0.81520800909731705
    "E2" = 0.20132728674488101;  "F2" = 0.85421984761270398;  "G2" = 0.133759827044514
    "H2" = 0.90388351848534998;  "I2" = 0.098501650956476894

    "B3" = 0.016170490914047202; "C3" = 0.0103050474808142;   "D3" = 0.89058884085795698
    "E3" = 0.10344453290089101;  "F3" = 0.89640639097119501;  "G3" = 0.096751885639587004
    "H3" = 0.93212189933187395;  "I3" = 0.077350965458864995

    "B4" = 0.0155172646863321;   "C4" = 0.013925957351606099; "D4" = 0.88662518804805002
    "E4" = 0.116588223474879;    "F4" = 0.89639866201201601;  "G4" = 0.10298873410117899
    "H4" = 0.94037570987261998;  "I4" = 0.073761035518276202

    "B5" = 0.0126835327477658;   "C5" = 0;                    "D5" = 0.94909090909090899
    "E5" = 0;                    "F5" = 0.96;                 "G5" = 0
    "H5" = 0.97955549955590704;  "I5" = 0
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

# C1 grew from a tiny fraction to a huge number this round -> display it in
# scientific notation so it still reads nicely next to the other columns.
$ws.Range("C1").NumberFormat = "0.00E+00"

# --- move / resize the first scatter chart ("Gráfico 1") -----------------
$chartObj = $ws.ChartObjects(1)
$chartObj.Left = 279.5
$chartObj.Top = 59.25
$chartObj.Width = 683.25
$chartObj.Height = 216

# --- selection now spans the whole refreshed table ------------------------
$ws.Range("A1:I5").Select()
